# Fix: Smart Signature Status. Return is_signed=False if tasks completed AFTER signature.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date (fecha) one day earlier
$ws.Range("A2").Value = 46021

# Update ticket_id to a simple numeric id instead of the huge float
$ws.Range("B2").Value = 123

# Move the active selection to I10 (matches saved cursor position)
$ws.Range("I10").Select()
